$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 111780628
$ws.Range("B2").Value = 78604
$ws.Range("D2").Value = "LC"
$ws.Range("F2").Value = "Norrlandslav"
$ws.Range("G2").Value = "Nephroma arcticum"
$ws.Range("H2").Value = "(L.) Torss."
$ws.Range("E2").Value = 6461
$ws.Range("Q2").Value = 707614
$ws.Range("R2").Value = 7397255
$ws.Range("Y2").Value = "'2023-08-29"
$ws.Range("Y2").ClearFormats()
$ws.Range("AA2").Value = "'2023-08-29"
$ws.Range("AA2").ClearFormats()
$ws.Range("Z2").ClearContents()
$ws.Range("AB2").ClearContents()

# Row 3
$ws.Range("A3").Value = 111780621
$ws.Range("B3").Value = 56543
$ws.Range("D3").Value = "NT"
$ws.Range("F3").Value = "Talltita"
$ws.Range("G3").Value = "Poecile montanus"
$ws.Range("H3").Value = "(Conrad von Baldenstein, 1827)"
$ws.Range("E3").Value = 103021
$ws.Range("Q3").Value = 707631
$ws.Range("R3").Value = 7397278
$ws.Range("Y3").Value = "'2023-08-29"
$ws.Range("Y3").ClearFormats()
$ws.Range("AA3").Value = "'2023-08-29"
$ws.Range("AA3").ClearFormats()
$ws.Range("Z3").ClearContents()
$ws.Range("AB3").ClearContents()

# Row 4
$ws.Range("A4").Value = 111780627
$ws.Range("B4").Value = 78604
$ws.Range("D4").Value = "LC"
$ws.Range("F4").Value = "Norrlandslav"
$ws.Range("G4").Value = "Nephroma arcticum"
$ws.Range("H4").Value = "(L.) Torss."
$ws.Range("E4").Value = 6461
$ws.Range("Q4").Value = 707647
$ws.Range("R4").Value = 7397287
$ws.Range("Y4").Value = "'2023-08-29"
$ws.Range("Y4").ClearFormats()
$ws.Range("AA4").Value = "'2023-08-29"
$ws.Range("AA4").ClearFormats()
$ws.Range("Z4").ClearContents()
$ws.Range("AB4").ClearContents()

# Row 5
$ws.Range("A5").Value = 111780624
$ws.Range("B5").Value = 95532
$ws.Range("D5").Value = "LC"
$ws.Range("F5").Value = "Revlummer"
$ws.Range("G5").Value = "Lycopodium annotinum"
$ws.Range("H5").Value = "L."
$ws.Range("E5").Value = 221945
$ws.Range("Q5").Value = 707601
$ws.Range("R5").Value = 7397313
$ws.Range("Y5").Value = "'2023-08-29"
$ws.Range("Y5").ClearFormats()
$ws.Range("AA5").Value = "'2023-08-29"
$ws.Range("AA5").ClearFormats()
$ws.Range("Z5").ClearContents()
$ws.Range("AB5").ClearContents()

# Row 6
$ws.Range("A6").Value = 111816118
$ws.Range("B6").Value = 78107
$ws.Range("D6").Value = "NT"
$ws.Range("F6").Value = "Vedskivlav"
$ws.Range("G6").Value = "Hertelidea botryosa"
$ws.Range("H6").Value = "(Fr.) Printzen & Kantvilas"
$ws.Range("E6").Value = 6453
$ws.Range("Q6").Value = 707670
$ws.Range("R6").Value = 7397328
$ws.Range("Y6").Value = "'2023-08-22"
$ws.Range("Y6").ClearFormats()
$ws.Range("AA6").Value = "'2023-08-22"
$ws.Range("AA6").ClearFormats()
$ws.Range("Z6").ClearContents()
$ws.Range("AB6").ClearContents()

# Row 7
$ws.Range("A7").Value = 111816132
$ws.Range("B7").Value = 95532
$ws.Range("D7").Value = "LC"
$ws.Range("F7").Value = "Revlummer"
$ws.Range("G7").Value = "Lycopodium annotinum"
$ws.Range("H7").Value = "L."
$ws.Range("E7").Value = 221945
$ws.Range("Q7").Value = 707590
$ws.Range("R7").Value = 7397240
$ws.Range("Y7").Value = "'2023-08-22"
$ws.Range("Y7").ClearFormats()
$ws.Range("AA7").Value = "'2023-08-22"
$ws.Range("AA7").ClearFormats()
$ws.Range("Z7").ClearContents()
$ws.Range("AB7").ClearContents()

# Row 8
$ws.Range("A8").Value = 111816145
$ws.Range("B8").Value = 77597
$ws.Range("D8").Value = "NT"
$ws.Range("F8").Value = "Knottrig blåslav"
$ws.Range("G8").Value = "Hypogymnia bitteri"
$ws.Range("H8").Value = "(Lynge) Ahti"
$ws.Range("E8").Value = 864
$ws.Range("Q8").Value = 707627
$ws.Range("R8").Value = 7397312
$ws.Range("Y8").Value = "'2023-08-22"
$ws.Range("Y8").ClearFormats()
$ws.Range("AA8").Value = "'2023-08-22"
$ws.Range("AA8").ClearFormats()
$ws.Range("Z8").ClearContents()
$ws.Range("AB8").ClearContents()

# Row 9
$ws.Range("A9").Value = 111816142
$ws.Range("B9").Value = 78604
$ws.Range("D9").Value = "LC"
$ws.Range("F9").Value = "Norrlandslav"
$ws.Range("G9").Value = "Nephroma arcticum"
$ws.Range("H9").Value = "(L.) Torss."
$ws.Range("E9").Value = 6461
$ws.Range("Q9").Value = 707613
$ws.Range("R9").Value = 7397270
$ws.Range("Y9").Value = "'2023-08-29"
$ws.Range("Y9").ClearFormats()
$ws.Range("AA9").Value = "'2023-08-29"
$ws.Range("AA9").ClearFormats()
$ws.Range("Z9").ClearContents()
$ws.Range("AB9").ClearContents()

# Row 10
$ws.Range("A10").Value = 111816137
$ws.Range("B10").Value = 90658
$ws.Range("D10").Value = "NT"
$ws.Range("F10").Value = "Orange taggsvamp"
$ws.Range("G10").Value = "Hydnellum aurantiacum"
$ws.Range("H10").Value = "(Batsch:Fr.) P.Karst."
$ws.Range("E10").Value = 4361
$ws.Range("Q10").Value = 707609
$ws.Range("R10").Value = 7397264
$ws.Range("Y10").Value = "'2023-08-22"
$ws.Range("Y10").ClearFormats()
$ws.Range("AA10").Value = "'2023-08-22"
$ws.Range("AA10").ClearFormats()
$ws.Range("Z10").ClearContents()
$ws.Range("AB10").ClearContents()

# Row 11
$ws.Range("A11").Value = 111816119
$ws.Range("B11").Value = 56543
$ws.Range("D11").Value = "NT"
$ws.Range("F11").Value = "Talltita"
$ws.Range("G11").Value = "Poecile montanus"
$ws.Range("H11").Value = "(Conrad von Baldenstein, 1827)"
$ws.Range("E11").Value = 103021
$ws.Range("Q11").Value = 707596
$ws.Range("R11").Value = 7397263
$ws.Range("Y11").Value = "'2023-08-22"
$ws.Range("Y11").ClearFormats()
$ws.Range("AA11").Value = "'2023-08-22"
$ws.Range("AA11").ClearFormats()
$ws.Range("Z11").ClearContents()
$ws.Range("AB11").ClearContents()

# Row 12
$ws.Range("A12").Value = 112202299
$ws.Range("B12").Value = 55611
$ws.Range("D12").Value = "NT"
$ws.Range("F12").Value = "Järpe"
$ws.Range("G12").Value = "Tetrastes bonasia"
$ws.Range("H12").Value = "(Linnaeus, 1758)"
$ws.Range("E12").Value = 102612
$ws.Range("Q12").Value = 707646
$ws.Range("R12").Value = 7397379
$ws.Range("Y12").Value = "'2023-09-17"
$ws.Range("Y12").ClearFormats()
$ws.Range("AA12").Value = "'2023-09-17"
$ws.Range("AA12").ClearFormats()
$ws.Range("Z12").ClearContents()
$ws.Range("AB12").ClearContents()
